$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 78.14530633333334
$ws.Range("H2").Value = 234.435919
$ws.Range("I2").Value = 0.9738103308619316
$ws.Range("J2").Value = 0.9738103308619316
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.321929333333333
$ws.Range("N2").Value = 21.965788
$ws.Range("Q2").Value = 572.1744107043525
$ws.Range("R2").Value = 5149.569696339172
$ws.Range("S2").Value = 0.9738103308619316
$ws.Range("T2").Value = 0.9738103308619316

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.9293790000000001
$ws.Range("H3").Value = 2.788137
$ws.Range("I3").Value = 0.01158148728249443
$ws.Range("J3").Value = 0.01158148728249443
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.321929333333333
$ws.Range("N3").Value = 21.965788
$ws.Range("Q3").Value = 6.804847361884001
$ws.Range("R3").Value = 61.24362625695601
$ws.Range("S3").Value = 0.01158148728249443
$ws.Range("T3").Value = 0.01158148728249443

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.172262
$ws.Range("H4").Value = 3.516786
$ws.Range("I4").Value = 0.01460818185557397
$ws.Range("J4").Value = 0.01460818185557397
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.321929333333333
$ws.Range("N4").Value = 21.965788
$ws.Range("Q4").Value = 8.583219524152002
$ws.Range("R4").Value = 77.248975717368
$ws.Range("S4").Value = 0.01460818185557397
$ws.Range("T4").Value = 0.01460818185557397
